# Slide 2's caption textbox ("TextBox 3") read a single run "Searching".
# Retitle it to "Search Party", split across three runs ("Search", " ",
# "Party") to match the author's edit.
$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$tr = $s.Shapes.Item(2).TextFrame.TextRange

$tr.Text = "Search"
[void]$tr.InsertAfter(" ")
[void]$tr.InsertAfter("Party")
